$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 41; this shifts existing rows 41-53 down to 42-54.
$ws.Rows.Item(41).Insert()

# Fill in the new row 41 with the weekly record (same static data as the
# other rows in this series, with its own date/volume/price figures).
$ws.Cells.Item(41, 1).Value = 4
$ws.Cells.Item(41, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(41, 3).Value = "Los Lagos"
$ws.Cells.Item(41, 4).Value = 44855
$ws.Cells.Item(41, 4).NumberFormat = $ws.Cells.Item(42, 4).NumberFormat
$ws.Cells.Item(41, 5).Value = 10
$ws.Cells.Item(41, 6).Value = 100112012
$ws.Cells.Item(41, 7).Value = "Espinaca"
$ws.Cells.Item(41, 8).Value = "Sin especificar"
$ws.Cells.Item(41, 9).Value = "Primera"
$ws.Cells.Item(41, 10).Value = 35
$ws.Cells.Item(41, 11).Value = 14000
$ws.Cells.Item(41, 12).Value = 14000
$ws.Cells.Item(41, 13).Value = 14000
$ws.Cells.Item(41, 14).Value = "$/cuna 10 kilos"
$ws.Cells.Item(41, 15).Value = "Región Metropolitana"
$ws.Cells.Item(41, 16).Value = 1400
$ws.Cells.Item(41, 17).Value = 10
$ws.Cells.Item(41, 18).Value = "Hortaliza"
